$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$apos = [char]39

# Row 2: replace 서귀포점 / 박광균 branch data with 양재점 / 점주 data
$ws.Range("A2").Value = "양재점"
$ws.Range("B2").Value = "점주"
$ws.Range("C2").Value = "점주(양재)"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "01099999999"
$ws.Range("E2").Value = "yangje@example.com"
$ws.Range("F2").Value = "휴대폰번호 중복,이메일 중복"

# Row 3: replace 남악점 / 정병훈 branch data with 양재점 / 사원 data (missing phone/email)
$ws.Range("A3").Value = "양재점"
$ws.Range("B3").Value = "사원"
$ws.Range("C3").Value = "사원(양재)"
$ws.Range("D3").Value = "$apos"
$ws.Range("E3").Value = "$apos"
$ws.Range("F3").Value = "필수입력값 누락,잘못된 이메일 형식,잘못된 휴대폰번호 형식"

# Rows 4-22 previously held additional branch rows; clear their contents
$ws.Range("A4:F22").ClearContents()
